$wb = $excel.ActiveWorkbook

# --- Update the daily conversion text on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.55 = 25845.35 pesos`n✅ 25845.35 pesos = 6.51 = 973.69 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 152.6
$ws2.Range("O10").Value = 3944
$ws2.Range("N12").Value = 3968.3
$ws2.Range("O12").Value = 149.5
